$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "note" (grade) column C with the actual marks awarded for
# each requirement in the correction grid.
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("C5").Value = 10
$ws.Range("C6").Value = 10
$ws.Range("C7").Value = 10
$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 10
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("C13").Value = 10
$ws.Range("C14").Value = 10
$ws.Range("C15").Value = 5

# Add a comment explaining the point deduction on the "plan du client
# suivie" row.
$ws.Range("D8").Value = "aucune section pour des jobs, plan du client non-remplis"

# Leave the selection on D9, matching where the grader was last working.
$ws.Range("D9").Select()
